# Add new columns I ("I0") and J ("IF") to the worksheet, matching the
# style of the existing header cells, and fill in the data for rows 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing header cell (H1) into the two new
# header cells so they pick up the same bold/border/centered style (s="1").
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

# Set the header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns, one triple (row, column, value) per new cell.
$newData = @(
    @(2, 9, 8),
    @(2, 10, 8),
    @(3, 9, 4),
    @(3, 10, 8),
    @(4, 9, 5),
    @(4, 10, 7),
    @(5, 9, 8),
    @(5, 10, 9),
    @(6, 9, 7),
    @(6, 10, 8),
    @(7, 9, 9),
    @(7, 10, 9),
    @(8, 9, 8),
    @(8, 10, 9),
    @(9, 9, 7),
    @(9, 10, 8),
    @(10, 9, 8),
    @(10, 10, 9),
    @(11, 9, 9),
    @(11, 10, 9),
    @(12, 9, 6),
    @(12, 10, 7),
    @(13, 9, 7),
    @(13, 10, 8),
    @(14, 9, 5),
    @(14, 10, 7),
    @(15, 9, 6),
    @(15, 10, 9),
    @(16, 9, 6),
    @(16, 10, 8),
    @(17, 9, 9),
    @(17, 10, 9),
    @(18, 9, 8),
    @(18, 10, 8),
    @(19, 9, 9),
    @(19, 10, 9),
    @(20, 9, 6),
    @(20, 10, 7),
    @(21, 9, 8),
    @(21, 10, 9),
    @(22, 9, 1),
    @(22, 10, 3),
    @(23, 9, 5),
    @(23, 10, 7),
    @(24, 9, 9),
    @(24, 10, 9),
    @(25, 9, 4),
    @(25, 10, 4),
    @(26, 9, 5),
    @(26, 10, 5),
    @(27, 9, 5),
    @(27, 10, 5),
    @(28, 9, 8),
    @(28, 10, 8),
    @(29, 9, 4),
    @(29, 10, 5),
    @(30, 9, 7),
    @(30, 10, 7),
    @(31, 9, 8),
    @(31, 10, 8),
    @(32, 9, 5),
    @(32, 10, 5),
    @(33, 9, 9),
    @(33, 10, 9),
    @(34, 9, 5),
    @(34, 10, 5),
)

foreach ($item in $newData) {
    $ws.Cells.Item($item[0], $item[1]).Value = $item[2]
}

Write-Host "Added columns I0 and J (IF) with $($newData.Count) cell updates"
